# Adds two new weekly price records to the "Locoto" price-report sheet.
# The table (header in row 1, data starting row 2) gets two brand new rows
# inserted: one before the current row 76, and one before the current row 87
# (i.e. right after the row that will have shifted into position 86).
# All existing rows below each insertion point shift down, which Excel does
# automatically (and preserves formatting, e.g. the date style on column D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First new record: inserted at row 76 ------------------------------
$ws.Rows.Item(76).Insert()

$ws.Range("A76").Value = 1
$ws.Range("B76").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C76").Value = "Arica y Parinacota"
$ws.Range("D76").Value = 44659
$ws.Range("E76").Value = 15
$ws.Range("F76").Value = 100112042
$ws.Range("G76").Value = "Locoto"
$ws.Range("H76").Value = "Sin especificar"
$ws.Range("I76").Value = "Segunda"
$ws.Range("J76").Value = 120
$ws.Range("K76").Value = 60000
$ws.Range("L76").Value = 62000
$ws.Range("M76").Value = 61000
$ws.Range("N76").Value = "`$/caja 20 kilos"
$ws.Range("O76").Value = "Región de Arica y Parinacota"
$ws.Range("P76").Value = 3050
$ws.Range("Q76").Value = 20
$ws.Range("R76").Value = "Hortaliza"

# --- Second new record: inserted at row 87 (after the first shift) -----
$ws.Rows.Item(87).Insert()

$ws.Range("A87").Value = 1
$ws.Range("B87").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C87").Value = "Arica y Parinacota"
$ws.Range("D87").Value = 44687
$ws.Range("E87").Value = 15
$ws.Range("F87").Value = 100112042
$ws.Range("G87").Value = "Locoto"
$ws.Range("H87").Value = "Sin especificar"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 100
$ws.Range("K87").Value = 58000
$ws.Range("L87").Value = 60000
$ws.Range("M87").Value = 59000
$ws.Range("N87").Value = "`$/caja 20 kilos"
$ws.Range("O87").Value = "Región de Arica y Parinacota"
$ws.Range("P87").Value = 2950
$ws.Range("Q87").Value = 20
$ws.Range("R87").Value = "Hortaliza"
